$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: after "...incremento de melhorias. " append a new sentence about
# the articles being versioned on GitHub, re-using the run immediately
# preceding it as a formatting template (Arial / 24 / cs-Arial / szCs-24).
# ---------------------------------------------------------------------------
$tmpl1 = $d.Content
$tmpl1.Find.Execute("e, com isso, permitir revisões, retrospectivas, registro de aprendizado e incremento de melhorias. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tmplStart1 = $tmpl1.Start
$tmplEnd1 = $tmpl1.End
$tmplLen1 = $tmplEnd1 - $tmplStart1

$insPos1 = $tmplEnd1
$dup1 = $d.Range($tmplStart1, $tmplEnd1)
$dup1.Copy()
$ins1 = $d.Range($insPos1, $insPos1)
$ins1.Paste()

$newRng1 = $d.Range($insPos1, $insPos1 + $tmplLen1)
$newRng1.Text = "E os artigos estarão versionados no GitHub para iteração."

# ---------------------------------------------------------------------------
# Edit 2: merge the "Passo " run and the "5" run into a single "Passo 5" run
# (the trailing ":" run is left alone). A same-text assignment is a no-op in
# this engine, so first rewrite to a distinct placeholder, then rewrite to
# the final text to force the underlying runs to coalesce.
# ---------------------------------------------------------------------------
$rngPasso5 = $d.Content
$rngPasso5.Find.Execute("Passo 5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p5Start = $rngPasso5.Start
$placeholderText = "Passo 5#placeholder#"
$rngPasso5.Text = $placeholderText
$rngPasso5b = $d.Range($p5Start, $p5Start + $placeholderText.Length)
$rngPasso5b.Text = "Passo 5"

# ---------------------------------------------------------------------------
# Edit 3: remove the first two trailing empty paragraphs (one list-style
# "PargrafodaLista" numbered item, one plain spacing paragraph), then fill
# the paragraph that follows (now the document's penultimate paragraph)
# with the new "Visão do Projeto" text, using a well-formatted run as the
# template again.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(39).Range.Delete()
$d.Paragraphs.Item(39).Range.Delete()

$tmpl3 = $d.Content
$tmpl3.Find.Execute("Retrospectiva;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tmplStart3 = $tmpl3.Start
$tmplEnd3 = $tmpl3.End
$tmplLen3 = $tmplEnd3 - $tmplStart3

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insPos3 = $lastPara.Range.Start

$dup3 = $d.Range($tmplStart3, $tmplEnd3)
$dup3.Copy()
$ins3 = $d.Range($insPos3, $insPos3)
$ins3.Paste()

$newRng3 = $d.Range($insPos3, $insPos3 + $tmplLen3)
$newText3 = "Passo 1 – A Visão do Projeto servirá para definir o objetivo do produto. É uma declaração capaz de explicar quem é o cliente, como o produto beneficiará esse cliente, quais funcionalidades serão essenciais e qual o diferencial competitivo frente ao mercado. E existem várias técnicas para sua elaboração, dentre as quais, serão utilizadas duas: Persona e Business Model Canvas."
$newRng3.Text = $newText3

Write-Output "ok"
